# Add data for 2021-11-25: roll the "through 11-16" snapshot forward to
# "through 11-17" and bump the November / Total figures for 2016-2021.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Rename the sheet and update its title text
$ws.Name = "Through 2021-11-17"

# Update the "November (through 11-16)" label in column A, row 12
$ws.Range("A12").Value = "November (through 11-17)"

# Update November row (row 12) for years 2016-2021 (columns C-H)
$ws.Range("C12").Value = 41
$ws.Range("D12").Value = 71
$ws.Range("E12").Value = 37
$ws.Range("F12").Value = 27
$ws.Range("G12").Value = 100
$ws.Range("H12").Value = 119

# Update Total row (row 13) for years 2016-2021 (columns C-H)
$ws.Range("C13").Value = 527
$ws.Range("D13").Value = 781
$ws.Range("E13").Value = 652
$ws.Range("F13").Value = 509
$ws.Range("G13").Value = 1157
$ws.Range("H13").Value = 1561
